$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("Z1")

function Set-TextValue($cellAddr, $text) {
    $helper.Formula = '="' + $text + '"'
    $helper.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
}

Set-TextValue "D2" "23.831.25"
Set-TextValue "E2" "  +2.17%  "
Set-TextValue "D3" "1.656.67"
Set-TextValue "E3" "  +1.94%  "
Set-TextValue "D4" "0.9981"
Set-TextValue "E4" "  -0.31%  "
Set-TextValue "D5" "0.9988"
Set-TextValue "E5" "  -0.21%  "
Set-TextValue "D6" "304.52"
Set-TextValue "E6" "  +0.54%  "
Set-TextValue "D7" "0.3830"
Set-TextValue "E7" "  +2.42%  "
Set-TextValue "D8" "0.3617"
Set-TextValue "E8" "  -0.09%  "
Set-TextValue "D9" "51.34"
Set-TextValue "E9" "  +0.15%  "
Set-TextValue "D10" "1.254"
Set-TextValue "E10" "  +2.39%  "
Set-TextValue "D11" "0.08245"
Set-TextValue "E11" "  +1.33%  "
Set-TextValue "D12" "0.9980"
Set-TextValue "E12" "  -0.31%  "
Set-TextValue "D13" "22.76"
Set-TextValue "E13" "  +2.08%  "
Set-TextValue "D14" "6.556"
Set-TextValue "E14" "  +1.25%  "
Set-TextValue "D15" "7.426"
Set-TextValue "E15" "  +2.11%  "
Set-TextValue "E16" "  +0.07%  "
Set-TextValue "D17" "1.653.95"
Set-TextValue "E17" "  +2.00%  "
Set-TextValue "D18" "97.61"
Set-TextValue "E18" "  +3.93%  "
Set-TextValue "D19" "0.06983"
Set-TextValue "E19" "  +0.68%  "
Set-TextValue "D20" "6.801"
Set-TextValue "E20" "  +4.29%  "
Set-TextValue "D21" "17.76"
Set-TextValue "E21" "  +1.43%  "
Set-TextValue "D22" "0.9985"
Set-TextValue "E22" "  -0.25%  "
Set-TextValue "D23" "12.69"
Set-TextValue "E23" "  +1.45%  "
Set-TextValue "D24" "23.827.80"
Set-TextValue "E24" "  +2.13%  "
Set-TextValue "D25" "2.567"
Set-TextValue "E25" "  +4.23%  "
Set-TextValue "D26" "3.122"
Set-TextValue "E26" "  +0.56%  "
Set-TextValue "D27" "21.33"
Set-TextValue "E27" "  +0.78%  "
Set-TextValue "D28" "151.05"
Set-TextValue "E28" "  +0.34%  "
Set-TextValue "D29" "5.238"
Set-TextValue "E29" "  -0.48%  "
Set-TextValue "D30" "135.01"
Set-TextValue "E30" "  +1.67%  "
Set-TextValue "D31" "1.838.85"
Set-TextValue "E31" "  +2.16%  "
Set-TextValue "D32" "7.021"
Set-TextValue "E32" "  +4.45%  "
Set-TextValue "D33" "1.083"
Set-TextValue "E33" "  +1.91%  "
Set-TextValue "D34" "11.99"
Set-TextValue "E34" "  +6.54%  "
Set-TextValue "D35" "2.118"
Set-TextValue "E35" "  -2.49%  "
Set-TextValue "D36" "0.02856"
Set-TextValue "E36" "  +3.57%  "
Set-TextValue "D37" "6.182"
Set-TextValue "E37" "  +3.09%  "
Set-TextValue "D38" "0.2524"
Set-TextValue "E38" "  +1.90%  "
Set-TextValue "D39" "0.08838"
Set-TextValue "E39" "  +0.94%  "
Set-TextValue "D40" "0.07091"
Set-TextValue "E40" "  -0.19%  "
Set-TextValue "D41" "12.89"
Set-TextValue "E41" "  +6.85%  "
Set-TextValue "D42" "0.7117"
Set-TextValue "E42" "  +2.04%  "
Set-TextValue "E43" "  +0.34%  "
Set-TextValue "D44" "15.87"
Set-TextValue "E44" "  -1.35%  "
Set-TextValue "D45" "0.6584"
Set-TextValue "E45" "  +1.81%  "
Set-TextValue "D46" "2.343"
Set-TextValue "E46" "  +3.28%  "
Set-TextValue "D47" "0.9984"
Set-TextValue "E47" "  -0.19%  "
Set-TextValue "D48" "3.980"
Set-TextValue "E48" "  +0.56%  "
Set-TextValue "D49" "0.07999"
Set-TextValue "E49" "  +0.38%  "
Set-TextValue "D50" "128.88"
Set-TextValue "E50" "  +2.43%  "
Set-TextValue "D51" "1.201"
Set-TextValue "E51" "  +1.52%  "

$helper.Clear()
$excel.CutCopyMode = $false
